$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.18"
$ws.Range("E2").Value = "'-0.50%"
$ws.Range("D3").Value = "'29.77"
$ws.Range("E3").Value = "'-1.26%"
$ws.Range("D4").Value = "'5.153"
$ws.Range("E4").Value = "'-0.49%"
$ws.Range("D5").Value = "'0.05765"
$ws.Range("E5").Value = "'0.31%"
$ws.Range("D6").Value = "'6.653"
$ws.Range("D7").Value = "'3.240"
$ws.Range("E7").Value = "'6.61%"
$ws.Range("D8").Value = "'0.8490"
$ws.Range("E8").Value = "'-1.18%"
$ws.Range("D9").Value = "'0.8533"
$ws.Range("E9").Value = "'-2.38%"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.01022"
$ws.Range("E10").Value = "'1,599.28%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1380"
$ws.Range("E11").Value = "'1.43%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07082"
$ws.Range("E12").Value = "'1.21%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03247"
$ws.Range("E13").Value = "'11.18%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09369"
$ws.Range("E14").Value = "'-0.33%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001540"
$ws.Range("E15").Value = "'1.15%"
$ws.Range("D16").Value = "'0.005968"
$ws.Range("E16").Value = "'-1.69%"
$ws.Range("D17").Value = "'3.521"
$ws.Range("E17").Value = "'0.45%"
$ws.Range("E18").Value = "'-1.77%"
$ws.Range("D19").Value = "'0.3164"
$ws.Range("E19").Value = "'-0.66%"
$ws.Range("D20").Value = "'0.03370"
$ws.Range("E20").Value = "'2.58%"
$ws.Range("E21").Value = "'0.72%"
$ws.Range("D22").Value = "'3.498"
$ws.Range("E22").Value = "'-3.16%"
$ws.Range("D23").Value = "'0.04136"
$ws.Range("E23").Value = "'0.16%"
$ws.Range("E24").Value = "'2.23%"
$ws.Range("D25").Value = "'0.001228"
$ws.Range("E25").Value = "'1.24%"
$ws.Range("E26").Value = "'-8.03%"
$ws.Range("D27").Value = "'0.0001199"
$ws.Range("E27").Value = "'1.73%"
$ws.Range("E28").Value = "'4.19%"
$ws.Range("D40").Value = "'0.03748"
$ws.Range("E40").Value = "'-1.04%"
$ws.Range("D41").Value = "'0.1070"
$ws.Range("E41").Value = "'-0.10%"
$ws.Range("D42").Value = "'0.002299"
$ws.Range("E42").Value = "'0.04%"
$ws.Range("D43").Value = "'0.002949"
$ws.Range("E43").Value = "'-48.95%"
$ws.Range("D44").Value = "'0.008543"
$ws.Range("E44").Value = "'-16.00%"
$ws.Range("D45").Value = "'0.00005500"
$ws.Range("E45").Value = "'8.19%"
$ws.Range("E46").Value = "'0.04%"
$ws.Range("D48").Value = "'0.002229"
$ws.Range("E48").Value = "'-18.77%"
$ws.Range("E49").Value = "'0.04%"
$ws.Range("E50").Value = "'0.04%"
